$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1427.0869
$ws.Cells.Item(112, 9).Value = 1050
$ws.Cells.Item(112, 10).Value = 1463
$ws.Cells.Item(112, 11).Value = 3150
$ws.Cells.Item(112, 12).Value = 4389
$ws.Cells.Item(112, 13).Value = -2042
$ws.Cells.Item(112, 14).Value = -6605
$ws.Cells.Item(123, 8).Value = 21607.691
$ws.Cells.Item(123, 9).Value = 15555
$ws.Cells.Item(123, 10).Value = 22112.084
$ws.Cells.Item(123, 11).Value = 15555
$ws.Cells.Item(123, 12).Value = 22112.084
$ws.Cells.Item(123, 13).Value = -10655
$ws.Cells.Item(123, 14).Value = -31912.084
$ws.Cells.Item(137, 8).Value = 1178.3243
$ws.Cells.Item(137, 9).Value = 996.9394
$ws.Cells.Item(137, 10).Value = 2674.75
$ws.Cells.Item(137, 11).Value = 2990.8182
$ws.Cells.Item(137, 12).Value = 8024.25
$ws.Cells.Item(137, 13).Value = -440.8181999999997
$ws.Cells.Item(137, 14).Value = -13124.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4811.1016
$ws.Cells.Item(32, 9).Value = 4099.143
$ws.Cells.Item(32, 10).Value = 7878
$ws.Cells.Item(32, 11).Value = 4099.143
$ws.Cells.Item(32, 12).Value = 7878
$ws.Cells.Item(32, 13).Value = -3812.143
$ws.Cells.Item(32, 14).Value = -8452
$ws.Cells.Item(61, 8).Value = 6819.6665
$ws.Cells.Item(61, 9).Value = 7060.65
$ws.Cells.Item(61, 11).Value = 7060.65
$ws.Cells.Item(61, 13).Value = -6848.65
$ws.Cells.Item(76, 8).Value = 30744.363
$ws.Cells.Item(76, 10).Value = 30744.363
$ws.Cells.Item(76, 12).Value = 30744.363
$ws.Cells.Item(76, 14).Value = -31420.363
$ws.Cells.Item(79, 8).Value = 30744.363
$ws.Cells.Item(79, 10).Value = 30744.363
$ws.Cells.Item(79, 12).Value = 30744.363
$ws.Cells.Item(79, 14).Value = -33084.363
$ws.Cells.Item(122, 8).Value = 3666147.2
$ws.Cells.Item(122, 9).Value = 4275836
$ws.Cells.Item(122, 11).Value = 12827508
$ws.Cells.Item(122, 13).Value = -12825058
$ws.Cells.Item(132, 8).Value = 3365.1
$ws.Cells.Item(132, 9).Value = 1409.5294
$ws.Cells.Item(132, 11).Value = 4228.5882
$ws.Cells.Item(132, 13).Value = -1698.5882
$ws.Cells.Item(136, 8).Value = 6819.6665
$ws.Cells.Item(136, 9).Value = 7060.65
$ws.Cells.Item(136, 11).Value = 21181.95
$ws.Cells.Item(136, 13).Value = -18631.95
$ws.Cells.Item(137, 8).Value = 39385
$ws.Cells.Item(137, 10).Value = 39385
$ws.Cells.Item(137, 12).Value = 39385
$ws.Cells.Item(137, 14).Value = -49585

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 49280
$ws.Cells.Item(59, 10).Value = 49280
$ws.Cells.Item(59, 12).Value = 49280
$ws.Cells.Item(59, 14).Value = -50974
$ws.Cells.Item(81, 8).Value = 30600
$ws.Cells.Item(81, 10).Value = 30600
$ws.Cells.Item(81, 12).Value = 30600
$ws.Cells.Item(81, 14).Value = -32722
$ws.Cells.Item(84, 8).Value = 30600
$ws.Cells.Item(84, 10).Value = 30600
$ws.Cells.Item(84, 12).Value = 91800
$ws.Cells.Item(84, 14).Value = -102408
$ws.Cells.Item(102, 8).Value = 32273
$ws.Cells.Item(102, 9).Value = 32273
$ws.Cells.Item(102, 11).Value = 32273
$ws.Cells.Item(102, 13).Value = -29028
$ws.Cells.Item(122, 8).Value = 22497.5
$ws.Cells.Item(122, 10).Value = 22497.5
$ws.Cells.Item(122, 12).Value = 22497.5
$ws.Cells.Item(122, 14).Value = -32297.5
$ws.Cells.Item(125, 8).Value = 52780
$ws.Cells.Item(125, 10).Value = 52780
$ws.Cells.Item(125, 12).Value = 52780
$ws.Cells.Item(125, 14).Value = -62620
$ws.Cells.Item(127, 8).Value = 56390
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 56390
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 13).Value = 56390
$ws.Cells.Item(127, 14).Value = -66310
$ws.Cells.Item(131, 8).Value = 52780
$ws.Cells.Item(131, 10).Value = 52780
$ws.Cells.Item(131, 12).Value = 52780
$ws.Cells.Item(131, 14).Value = -62860

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1404.8334
$ws.Cells.Item(16, 9).Value = 1666.3334
$ws.Cells.Item(16, 10).Value = 1317.6666
$ws.Cells.Item(16, 11).Value = 1666.3334
$ws.Cells.Item(16, 12).Value = 1317.6666
$ws.Cells.Item(16, 13).Value = -1379.3334
$ws.Cells.Item(16, 14).Value = -1891.6666
$ws.Cells.Item(92, 8).Value = 39993.332
$ws.Cells.Item(92, 10).Value = 39993.332
$ws.Cells.Item(92, 12).Value = 39993.332
$ws.Cells.Item(92, 14).Value = -44985.332
$ws.Cells.Item(113, 8).Value = 1404.8334
$ws.Cells.Item(113, 9).Value = 1666.3334
$ws.Cells.Item(113, 10).Value = 1317.6666
$ws.Cells.Item(113, 11).Value = 1666.3334
$ws.Cells.Item(113, 12).Value = 1317.6666
$ws.Cells.Item(113, 13).Value = 503.6666
$ws.Cells.Item(113, 14).Value = -5657.6666
$ws.Cells.Item(132, 8).Value = 2343.1304
$ws.Cells.Item(132, 9).Value = 1995.75
$ws.Cells.Item(132, 10).Value = 4659
$ws.Cells.Item(132, 11).Value = 5987.25
$ws.Cells.Item(132, 12).Value = 13977
$ws.Cells.Item(132, 13).Value = -3457.25
$ws.Cells.Item(132, 14).Value = -19037
$ws.Cells.Item(135, 8).Value = 35462.35
$ws.Cells.Item(135, 10).Value = 35462.35
$ws.Cells.Item(135, 12).Value = 35462.35
$ws.Cells.Item(135, 14).Value = -45602.35

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 545827.25
$ws.Cells.Item(5, 9).Value = 800
$ws.Cells.Item(5, 10).Value = 750212.5
$ws.Cells.Item(5, 11).Value = 2400
$ws.Cells.Item(5, 12).Value = 2250637.5
$ws.Cells.Item(5, 13).Value = -2288
$ws.Cells.Item(5, 14).Value = -2250861.5
$ws.Cells.Item(68, 8).Value = 2344.4
$ws.Cells.Item(68, 9).Value = 430.5
$ws.Cells.Item(68, 11).Value = 1291.5
$ws.Cells.Item(68, 13).Value = -480.5
$ws.Cells.Item(71, 8).Value = 2344.4
$ws.Cells.Item(71, 9).Value = 430.5
$ws.Cells.Item(71, 11).Value = 3874.5
$ws.Cells.Item(71, 13).Value = 181.5
$ws.Cells.Item(107, 8).Value = 111822.28
$ws.Cells.Item(107, 9).Value = 200
$ws.Cells.Item(107, 10).Value = 125775.06
$ws.Cells.Item(107, 11).Value = 600
$ws.Cells.Item(107, 12).Value = 377325.18
$ws.Cells.Item(107, 13).Value = 1320
$ws.Cells.Item(107, 14).Value = -381165.18
$ws.Cells.Item(122, 8).Value = 786.2857
$ws.Cells.Item(122, 10).Value = 1005
$ws.Cells.Item(122, 12).Value = 9045
$ws.Cells.Item(122, 14).Value = -13945
$ws.Cells.Item(125, 8).Value = 3961.5334
$ws.Cells.Item(125, 9).Value = 1500
$ws.Cells.Item(125, 10).Value = 4137.357
$ws.Cells.Item(125, 11).Value = 4500
$ws.Cells.Item(125, 12).Value = 12412.071
$ws.Cells.Item(125, 13).Value = 420
$ws.Cells.Item(125, 14).Value = -22252.071
$ws.Cells.Item(131, 8).Value = 1516099.8
$ws.Cells.Item(131, 9).Value = 5882968.5
$ws.Cells.Item(131, 10).Value = 1063.6123
$ws.Cells.Item(131, 11).Value = 17648905.5
$ws.Cells.Item(131, 12).Value = 3190.8369
$ws.Cells.Item(131, 13).Value = -17643865.5
$ws.Cells.Item(131, 14).Value = -13270.8369
$ws.Cells.Item(135, 8).Value = 545827.25
$ws.Cells.Item(135, 9).Value = 800
$ws.Cells.Item(135, 10).Value = 750212.5
$ws.Cells.Item(135, 11).Value = 7200
$ws.Cells.Item(135, 12).Value = 6751912.5
$ws.Cells.Item(135, 13).Value = -4665
$ws.Cells.Item(135, 14).Value = -6756982.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3104950
$ws.Cells.Item(122, 9).Value = 2236857
$ws.Cells.Item(122, 10).Value = 6251787.5
$ws.Cells.Item(122, 11).Value = 6710571
$ws.Cells.Item(122, 12).Value = 18755362.5
$ws.Cells.Item(122, 13).Value = -6708121
$ws.Cells.Item(122, 14).Value = -18760262.5
$ws.Cells.Item(132, 8).Value = 3725.7585
$ws.Cells.Item(132, 9).Value = 4508.067
$ws.Cells.Item(132, 10).Value = 2887.5715
$ws.Cells.Item(132, 11).Value = 13524.201
$ws.Cells.Item(132, 12).Value = 8662.7145
$ws.Cells.Item(132, 13).Value = -10994.201
$ws.Cells.Item(132, 14).Value = -13722.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 12346585
$ws.Cells.Item(46, 9).Value = 18519210
$ws.Cells.Item(46, 10).Value = 1333.1111
$ws.Cells.Item(46, 11).Value = 18519210
$ws.Cells.Item(46, 12).Value = 1333.1111
$ws.Cells.Item(46, 13).Value = -18519022
$ws.Cells.Item(46, 14).Value = -1709.1111
$ws.Cells.Item(122, 8).Value = 11907179
$ws.Cells.Item(122, 10).Value = 2876.25
$ws.Cells.Item(122, 12).Value = 8628.75
$ws.Cells.Item(122, 14).Value = -13528.75
